{"js": "// Update the division-problem worksheet cells: each <w:t> in the table\n// holds a unique \"NNN\u00f7N=\" string. Replace each old problem with its new\n// counterpart per the commit diff, using Word's search/replace so the\n// surrounding run formatting (font, size) is preserved untouched.\nconst replacements = [\n  [\"663\u00f74=\", \"443\u00f78=\"],\n  [\"896\u00f78=\", \"971\u00f78=\"],\n  [\"716\u00f72=\", \"619\u00f79=\"],\n  [\"333\u00f75=\", \"779\u00f78=\"],\n  [\"665\u00f76=\", \"810\u00f76=\"],\n  [\"344\u00f77=\", \"230\u00f78=\"],\n  [\"439\u00f76=\", \"570\u00f79=\"],\n  [\"284\u00f78=\", \"514\u00f74=\"],\n  [\"161\u00f73=\", \"524\u00f75=\"],\n  [\"641\u00f78=\", \"400\u00f74=\"],\n  [\"518\u00f72=\", \"425\u00f74=\"],\n  [\"400\u00f72=\", \"361\u00f79=\"],\n  [\"582\u00f75=\", \"950\u00f78=\"],\n  [\"704\u00f76=\", \"676\u00f76=\"],\n  [\"722\u00f73=\", \"143\u00f77=\"],\n  [\"687\u00f76=\", \"546\u00f76=\"],\n  [\"557\u00f79=\", \"366\u00f74=\"],\n  [\"495\u00f73=\", \"861\u00f77=\"],\n  [\"507\u00f73=\", \"454\u00f78=\"],\n  [\"163\u00f78=\", \"824\u00f73=\"],\n  [\"297\u00f72=\", \"610\u00f75=\"],\n  [\"409\u00f76=\", \"838\u00f77=\"],\n  [\"626\u00f76=\", \"322\u00f78=\"],\n  [\"261\u00f72=\", \"362\u00f73=\"],\n  [\"750\u00f73=\", \"906\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem worksheet cells: each cell in the table\n# holds a unique \"NNN\u00f7N=\" run of text. Replace each old problem with its\n# new counterpart per the commit diff, using Word's Find/Replace so the\n# surrounding run formatting (font, size) is left untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"663\u00f74=\"; new=\"443\u00f78=\"},\n  @{old=\"896\u00f78=\"; new=\"971\u00f78=\"},\n  @{old=\"716\u00f72=\"; new=\"619\u00f79=\"},\n  @{old=\"333\u00f75=\"; new=\"779\u00f78=\"},\n  @{old=\"665\u00f76=\"; new=\"810\u00f76=\"},\n  @{old=\"344\u00f77=\"; new=\"230\u00f78=\"},\n  @{old=\"439\u00f76=\"; new=\"570\u00f79=\"},\n  @{old=\"284\u00f78=\"; new=\"514\u00f74=\"},\n  @{old=\"161\u00f73=\"; new=\"524\u00f75=\"},\n  @{old=\"641\u00f78=\"; new=\"400\u00f74=\"},\n  @{old=\"518\u00f72=\"; new=\"425\u00f74=\"},\n  @{old=\"400\u00f72=\"; new=\"361\u00f79=\"},\n  @{old=\"582\u00f75=\"; new=\"950\u00f78=\"},\n  @{old=\"704\u00f76=\"; new=\"676\u00f76=\"},\n  @{old=\"722\u00f73=\"; new=\"143\u00f77=\"},\n  @{old=\"687\u00f76=\"; new=\"546\u00f76=\"},\n  @{old=\"557\u00f79=\"; new=\"366\u00f74=\"},\n  @{old=\"495\u00f73=\"; new=\"861\u00f77=\"},\n  @{old=\"507\u00f73=\"; new=\"454\u00f78=\"},\n  @{old=\"163\u00f78=\"; new=\"824\u00f73=\"},\n  @{old=\"297\u00f72=\"; new=\"610\u00f75=\"},\n  @{old=\"409\u00f76=\"; new=\"838\u00f77=\"},\n  @{old=\"626\u00f76=\"; new=\"322\u00f78=\"},\n  @{old=\"261\u00f72=\"; new=\"362\u00f73=\"},\n  @{old=\"750\u00f73=\"; new=\"906\u00f72=\"}\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $p.old\n  $find.Replacement.Text = $p.new\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute([ref]$p.old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$p.new, [ref]2) | Out-Null\n}\n"}
